# Add newly-deleted "Direct Steam Tower Receiver" SSC variables to the
# "SAM Variable Changes" sheet (rows 49-56), mirroring the format already
# used for the "Direct Steam Tower Parasitics" variables in rows 43-48.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# Copy the per-column formatting from row 43 (an existing "Deleted
# variable" row) down onto the 8 new rows, column by column so we don't
# touch column D (which has no cell in these rows).
$formatCols = @("A", "B", "C", "E", "F", "G", "H")
foreach ($col in $formatCols) {
    $ws.Range("$col`43").Copy()
    $ws.Range("$col`49:$col`56").PasteSpecial(-4122)
}

# New deleted variables, all belonging to the "Direct Steam Tower
# Receiver" input page, same as the existing rows above them.
$deletedVars = @(
    "rec_htf_flow",
    "ref_htf_flow",
    "P_cold_tank",
    "P_tower_conv",
    "P_tower_rad",
    "P_htf_pump",
    "night_recirc",
    "mode"
)

for ($i = 0; $i -lt $deletedVars.Length; $i++) {
    $row = 49 + $i
    $ws.Range("A$row").Value = "Deleted variable"
    $ws.Range("B$row").Value = "number"
    $ws.Range("C$row").Value = $deletedVars[$i]
    $ws.Range("E$row").Value = "Direct Steam Tower Receiver"
    $ws.Range("F$row").Value = "not used"
    $ws.Range("G$row").Value = "N/A"
    $ws.Range("H$row").Value = "Ty"
}

# Match the author's final cursor/selection position recorded in the
# saved file.
$ws.Activate()
$ws.Range("A57").Select()
